$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...par rapport à un tirage." ->
#           "...par rapport à un tirage" + " donné" + "." (3 runs, same font)
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("tirage.", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'tirage.' target sentence"
}

# Range covering just "tirage" (drop the trailing period from the match).
$tirage = $d.Range($r.Start, $r.End - 1)
$tirage.Collapse(0)                 # wdCollapseEnd -> collapse to right after "tirage"
$tirage.InsertAfter(" donné")       # becomes part of the same run for now

# Range covering the newly-inserted " donné" text.
$donne = $d.Range($tirage.Start, $tirage.Start + 6)
# Toggle Bold on/off: this is a no-op in appearance but forces the engine to
# split " donné" into its own run (preserving the full original rFonts/size).
$donne.Bold = $true
$donne.Bold = $false

# The trailing "." is still sitting right after " donné" (original char,
# untouched so far) - split it into its own run the same way.
$period = $d.Range($donne.End, $donne.End + 1)
$period.Bold = $true
$period.Bold = $false

# ---------------------------------------------------------------------------
# Change 2: merge "apply plugin: 'kotlin-kapt'" + " " into a single run
#           "apply plugin: 'kotlin-kapt' " (trailing space, same run)
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("apply plugin: 'kotlin-kapt'", $false, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the gradle code snippet run"
}

# The single space that lives in its own run right after the code text.
$spaceRun = $d.Range($r2.End, $r2.End + 1)
$spaceRun.Delete()

# Re-add the space right after the code text - since the following run now
# carries identical (Arial) formatting, not the Consolas one, this extends
# the code run itself instead of merging into the Arial text.
$r2.InsertAfter(" ")
